$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: becomes the "VU / Knärot" record, with earlier time (18:29) ---
$ws.Range("A2").Value = 111790625
$ws.Range("B2").Value = 96348
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 220787
$ws.Range("F2").Value = "Knärot"
$ws.Range("G2").Value = "Goodyera repens"
$ws.Range("H2").Value = "(L.) R. Br."
$ws.Range("Q2").Value = 489824.6884970492
$ws.Range("R2").Value = 6949020.70113107
$ws.Range("Z2").Value = "18:29"
$ws.Range("AB2").Value = "18:29"
$ws.Range("AC2").ClearContents()

# --- Row 4: becomes the "NT / Garnlav" record, with later time (18:34) ---
$ws.Range("A4").Value = 111790785
$ws.Range("B4").Value = 77515
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 489818.2822038208
$ws.Range("R4").Value = 6949032.207674611
$ws.Range("Z4").Value = "18:34"
$ws.Range("AB4").Value = "18:34"
$ws.Range("AC4").Value = "Många träd med mycket lav i området"
